$wb = $excel.ActiveWorkbook

# --- Sheet: Daily (row 2) ---
$daily = $wb.Worksheets.Item("Daily")
$daily.Range("G2").Value = 2945.83
$daily.Range("H2").Value = 6212.94
$daily.Range("I2").Value = 726.63
$daily.Range("J2").Value = 736.54
$daily.Range("L2").Value = 736.54

# --- Sheet: Hourly (rows 9-19) ---
$hourly = $wb.Worksheets.Item("Hourly")

$hourly.Range("I9").Value = 30.52
$hourly.Range("K9").Value = 0.9
$hourly.Range("M9").Value = 0.9

$hourly.Range("H10").Value = 94.5
$hourly.Range("I10").Value = 402.86
$hourly.Range("J10").Value = 45.62
$hourly.Range("K10").Value = 23.62
$hourly.Range("M10").Value = 23.62

$hourly.Range("H11").Value = 238.15
$hourly.Range("I11").Value = 626.09
$hourly.Range("J11").Value = 70.65000000000001
$hourly.Range("K11").Value = 59.54
$hourly.Range("M11").Value = 59.54

$hourly.Range("H12").Value = 363.34
$hourly.Range("I12").Value = 729.99
$hourly.Range("J12").Value = 84.93000000000001
$hourly.Range("K12").Value = 90.83
$hourly.Range("M12").Value = 90.83

$hourly.Range("H13").Value = 449.21
$hourly.Range("I13").Value = 781.7
$hourly.Range("J13").Value = 92.88
$hourly.Range("K13").Value = 112.3
$hourly.Range("M13").Value = 112.3

$hourly.Range("H14").Value = 485.22
$hourly.Range("I14").Value = 800.37
$hourly.Range("J14").Value = 95.89
$hourly.Range("K14").Value = 121.31
$hourly.Range("M14").Value = 121.31

$hourly.Range("H15").Value = 467.33
$hourly.Range("I15").Value = 791.35
$hourly.Range("J15").Value = 94.40000000000001
$hourly.Range("K15").Value = 116.83
$hourly.Range("M15").Value = 116.83

$hourly.Range("H16").Value = 397.5
$hourly.Range("I16").Value = 752.09
$hourly.Range("J16").Value = 88.2
$hourly.Range("K16").Value = 99.38
$hourly.Range("M16").Value = 99.38

$hourly.Range("H17").Value = 283.99
$hourly.Range("I17").Value = 669.6
$hourly.Range("J17").Value = 76.33
$hourly.Range("K17").Value = 71
$hourly.Range("M17").Value = 71

$hourly.Range("H18").Value = 143.09
$hourly.Range("I18").Value = 501.66
$hourly.Range("J18").Value = 55.91

$hourly.Range("I19").Value = 126.71
$hourly.Range("K19").Value = 5.06
$hourly.Range("M19").Value = 5.06
